# "R vs Python 3, 4, 5" — add two new rows of R/Python equivalents
# (order/pd.sort_value and read.csv/pd.read_csv) right after the
# "install.packages/import" row, pushing the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R vs Python")

# Insert two blank rows above the current row 4 ("class"/"type"), shifting
# everything from row 4 down to row 6 onward.
$ws.Rows("4:5").Insert()

# New row 4: order / pd.sort_value
$ws.Range("A4").Value = "order"
$ws.Range("B4").Value = "pd.sort_value"
$ws.Range("C4").Value = "Done"
$ws.Range("D4").Value = "CX"

# New row 5: read.csv / pd.read_csv
$ws.Range("A5").Value = "read.csv"
$ws.Range("B5").Value = "pd.read_csv"
$ws.Range("C5").Value = "Done"
$ws.Range("D5").Value = "CX"

# Restore the sheet's active cell/selection as left by the author.
$ws.Range("D8").Select()
